$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark that used to sit at the end of
#        the last non-empty paragraph ("...1db futtatható ciklus"). Do
#        this FIRST so there is never a moment with two same-named
#        "_GoBack" bookmarks (which would make the later lookup
#        ambiguous). ---

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Append ", custom scroll" (with proofErr spell markers) and a
#        new zero-width "_GoBack" bookmark to the end of the FIRST
#        paragraph, right after " hogy melyik a jó kód 2 közül)". ---

$p1 = $d.Paragraphs(1)
$full = $p1.Range
$full.End = $full.End - 1      # exclude the paragraph mark
$full.Collapse(0)              # collapse to the very end of the text

# Use a one-character throwaway run as an anchor: InsertXML replaces the
# text of the exact range it is called on and re-appends new content at
# the end of the paragraph, so inserting a harmless placeholder char first
# (which leaves all the *existing* runs fully untouched) and then
# replacing *just that* placeholder lets us splice in new runs (incl.
# <w:proofErr/> and bookmark markers) without disturbing anything already
# in the paragraph.
$full.InsertAfter("X")
$placeholder = $d.Range($full.End - 1, $full.End)

$newRunsXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>custom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> scroll</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$placeholder.InsertXML($newRunsXml)
